$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.68463133333333
$ws.Range("H2").Value = 35.053894
$ws.Range("I2").Value = 0.1650073836014231
$ws.Range("J2").Value = 0.165007383601423
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 12.840326
$ws.Range("N2").Value = 38.520978
$ws.Range("O2").Value = 0.3393128690704512
$ws.Range("P2").Value = 0.3393128690704511
$ws.Range("Q2").Value = 150.0344755098147
$ws.Range("R2").Value = 1350.310279588332
$ws.Range("S2").Value = 0.05598912874760739
$ws.Range("T2").Value = 0.05598912874760735

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.68463133333333
$ws.Range("H3").Value = 35.053894
$ws.Range("I3").Value = 0.1650073836014231
$ws.Range("J3").Value = 0.165007383601423
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.834223333333334
$ws.Range("N3").Value = 29.50267
$ws.Range("O3").Value = 0.2598749077175229
$ws.Range("P3").Value = 0.2598749077175228
$ws.Range("Q3").Value = 114.9092740996645
$ws.Range("R3").Value = 1034.18346689698
$ws.Range("S3").Value = 0.04288127858612972
$ws.Range("T3").Value = 0.0428812785861297

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.68463133333333
$ws.Range("H4").Value = 35.053894
$ws.Range("I4").Value = 0.1650073836014231
$ws.Range("J4").Value = 0.165007383601423
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.654269333333334
$ws.Range("N4").Value = 25.962808
$ws.Range("O4").Value = 0.228693956617749
$ws.Range("P4").Value = 0.2286939566177489
$ws.Range("Q4").Value = 101.1219466193725
$ws.Range("R4").Value = 910.097519574352
$ws.Range("S4").Value = 0.03773619142695211
$ws.Range("T4").Value = 0.03773619142695209

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.68463133333333
$ws.Range("H5").Value = 35.053894
$ws.Range("I5").Value = 0.1650073836014231
$ws.Range("J5").Value = 0.165007383601423
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.513324000000001
$ws.Range("N5").Value = 19.539972
$ws.Range("O5").Value = 0.1721182665942771
$ws.Range("P5").Value = 0.1721182665942771
$ws.Range("Q5").Value = 76.10578969455202
$ws.Range("R5").Value = 684.9521072509681
$ws.Range("S5").Value = 0.02840078484073389
$ws.Range("T5").Value = 0.02840078484073387

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.82343466666667
$ws.Range("H6").Value = 47.470304
$ws.Range("I6").Value = 0.2234545087003506
$ws.Range("J6").Value = 0.2234545087003506
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 12.840326
$ws.Range("N6").Value = 38.520978
$ws.Range("O6").Value = 0.3393128690704512
$ws.Range("P6").Value = 0.3393128690704511
$ws.Range("Q6").Value = 203.1780595597013
$ws.Range("R6").Value = 1828.602536037312
$ws.Range("S6").Value = 0.07582099045384406
$ws.Range("T6").Value = 0.07582099045384402

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.82343466666667
$ws.Range("H7").Value = 47.470304
$ws.Range("I7").Value = 0.2234545087003506
$ws.Range("J7").Value = 0.2234545087003506
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.834223333333334
$ws.Range("N7").Value = 29.50267
$ws.Range("O7").Value = 0.2598749077175229
$ws.Range("P7").Value = 0.2598749077175228
$ws.Range("Q7").Value = 155.6111904124089
$ws.Range("R7").Value = 1400.50071371168
$ws.Range("S7").Value = 0.05807021982756801
$ws.Range("T7").Value = 0.058070219827568

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 15.82343466666667
$ws.Range("H8").Value = 47.470304
$ws.Range("I8").Value = 0.2234545087003506
$ws.Range("J8").Value = 0.2234545087003506
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 8.654269333333334
$ws.Range("N8").Value = 25.962808
$ws.Range("O8").Value = 0.228693956617749
$ws.Range("P8").Value = 0.2286939566177489
$ws.Range("Q8").Value = 136.9402653837369
$ws.Range("R8").Value = 1232.462388453632
$ws.Range("S8").Value = 0.05110269571875839
$ws.Range("T8").Value = 0.05110269571875837

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 15.82343466666667
$ws.Range("H9").Value = 47.470304
$ws.Range("I9").Value = 0.2234545087003506
$ws.Range("J9").Value = 0.2234545087003506
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.513324000000001
$ws.Range("N9").Value = 19.539972
$ws.Range("O9").Value = 0.1721182665942771
$ws.Range("P9").Value = 0.1721182665942771
$ws.Range("Q9").Value = 103.063156776832
$ws.Range("R9").Value = 927.5684109914881
$ws.Range("S9").Value = 0.03846060270018016
$ws.Range("T9").Value = 0.03846060270018014

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 36.36549133333333
$ws.Range("H10").Value = 109.096474
$ws.Range("I10").Value = 0.5135441938313808
$ws.Range("J10").Value = 0.5135441938313807
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 12.840326
$ws.Range("N10").Value = 38.520978
$ws.Range("O10").Value = 0.3393128690704512
$ws.Range("P10").Value = 0.3393128690704511
$ws.Range("Q10").Value = 466.9447638701746
$ws.Range("R10").Value = 4202.502874831572
$ws.Range("S10").Value = 0.1742521538033978
$ws.Range("T10").Value = 0.1742521538033976

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 36.36549133333333
$ws.Range("H11").Value = 109.096474
$ws.Range("I11").Value = 0.5135441938313808
$ws.Range("J11").Value = 0.5135441938313807
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 9.834223333333334
$ws.Range("N11").Value = 29.50267
$ws.Range("O11").Value = 0.2598749077175229
$ws.Range("P11").Value = 0.2598749077175228
$ws.Range("Q11").Value = 357.6263633983978
$ws.Range("R11").Value = 3218.63727058558
$ws.Range("S11").Value = 0.1334572499807998
$ws.Range("T11").Value = 0.1334572499807997

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 36.36549133333333
$ws.Range("H12").Value = 109.096474
$ws.Range("I12").Value = 0.5135441938313808
$ws.Range("J12").Value = 0.5135441938313807
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 8.654269333333334
$ws.Range("N12").Value = 25.962808
$ws.Range("O12").Value = 0.228693956617749
$ws.Range("P12").Value = 0.2286939566177489
$ws.Range("Q12").Value = 314.7167564376658
$ws.Range("R12").Value = 2832.450807938992
$ws.Range("S12").Value = 0.1174444535853707
$ws.Range("T12").Value = 0.1174444535853706

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 36.36549133333333
$ws.Range("H13").Value = 109.096474
$ws.Range("I13").Value = 0.5135441938313808
$ws.Range("J13").Value = 0.5135441938313807
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.513324000000001
$ws.Range("N13").Value = 19.539972
$ws.Range("O13").Value = 0.1721182665942771
$ws.Range("P13").Value = 0.1721182665942771
$ws.Range("Q13").Value = 236.860227473192
$ws.Range("R13").Value = 2131.742047258728
$ws.Range("S13").Value = 0.08839033646181273
$ws.Range("T13").Value = 0.08839033646181269

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.939221333333332
$ws.Range("H14").Value = 20.817664
$ws.Range("I14").Value = 0.09799391386684557
$ws.Range("J14").Value = 0.09799391386684556
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 12.840326
$ws.Range("N14").Value = 38.520978
$ws.Range("O14").Value = 0.3393128690704512
$ws.Range("P14").Value = 0.3393128690704511
$ws.Range("Q14").Value = 89.10186410615465
$ws.Range("R14").Value = 801.9167769553919
$ws.Range("S14").Value = 0.03325059606560204
$ws.Range("T14").Value = 0.03325059606560203

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.939221333333332
$ws.Range("H15").Value = 20.817664
$ws.Range("I15").Value = 0.09799391386684557
$ws.Range("J15").Value = 0.09799391386684556
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 9.834223333333334
$ws.Range("N15").Value = 29.50267
$ws.Range("O15").Value = 0.2598749077175229
$ws.Range("P15").Value = 0.2598749077175228
$ws.Range("Q15").Value = 68.2418523514311
$ws.Range("R15").Value = 614.17667116288
$ws.Range("S15").Value = 0.02546615932302538
$ws.Range("T15").Value = 0.02546615932302537

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.939221333333332
$ws.Range("H16").Value = 20.817664
$ws.Range("I16").Value = 0.09799391386684557
$ws.Range("J16").Value = 0.09799391386684556
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 8.654269333333334
$ws.Range("N16").Value = 25.962808
$ws.Range("O16").Value = 0.228693956617749
$ws.Range("P16").Value = 0.2286939566177489
$ws.Range("Q16").Value = 60.05389038227911
$ws.Range("R16").Value = 540.4850134405119
$ws.Range("S16").Value = 0.02241061588666781
$ws.Range("T16").Value = 0.0224106158866678

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.939221333333332
$ws.Range("H17").Value = 20.817664
$ws.Range("I17").Value = 0.09799391386684557
$ws.Range("J17").Value = 0.09799391386684556
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.513324000000001
$ws.Range("N17").Value = 19.539972
$ws.Range("O17").Value = 0.1721182665942771
$ws.Range("P17").Value = 0.1721182665942771
$ws.Range("Q17").Value = 45.197396851712
$ws.Range("R17").Value = 406.776571665408
$ws.Range("S17").Value = 0.01686654259155036
$ws.Range("T17").Value = 0.01686654259155035

Write-Output "applied edits"